$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): set F1, and add G1/H1 with header formatting copied from A1 ---
$ws.Range("F1").Value = "Årsag"
$ws.Range("G1").Value = "Ny leverandør"
$ws.Range("H1").Value = "TCV_range"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows 2-40: update column F (Arsag), set G (Ny leverandoer) where applicable, set H (TCV_range) ---
# Row 2
$ws.Range("F2").Value = "Ikke oplyst"
$ws.Range("H2").Value = "20000-40000"
# Row 3
$ws.Range("F3").Value = "Ikke flere medarbejdere i virksomheden"
$ws.Range("H3").Value = "20000-40000"
# Row 4
$ws.Range("F4").Value = "Pris"
$ws.Range("G4").Value = "DataLøn"
$ws.Range("H4").Value = "20000-40000"
# Row 5
$ws.Range("F5").Value = "Pris"
$ws.Range("G5").Value = "DataLøn"
$ws.Range("H5").Value = "20000-40000"
# Row 6
$ws.Range("F6").Value = "Ikke oplyst"
$ws.Range("H6").Value = "20000-40000"
# Row 7
$ws.Range("F7").Value = "Covid-19"
$ws.Range("H7").Value = "20000-40000"
# Row 8
$ws.Range("F8").Value = "Bruger ikke produktet"
$ws.Range("H8").Value = "20000-40000"
# Row 9
$ws.Range("F9").Value = "Ikke oplyst"
$ws.Range("H9").Value = "20000-40000"
# Row 10
$ws.Range("F10").Value = "Ikke oplyst"
$ws.Range("H10").Value = "20000-40000"
# Row 11
$ws.Range("F11").Value = "Ikke oplyst"
$ws.Range("H11").Value = "20000-40000"
# Row 12
$ws.Range("F12").Value = "Virksomheden lukker"
$ws.Range("H12").Value = "20000-40000"
# Row 13
$ws.Range("F13").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("G13").Value = "DataLøn"
$ws.Range("H13").Value = "20000-40000"
# Row 14
$ws.Range("F14").Value = "Ikke oplyst"
$ws.Range("H14").Value = "20000-40000"
# Row 15
$ws.Range("F15").Value = "Ikke oplyst"
$ws.Range("H15").Value = "20000-40000"
# Row 16
$ws.Range("F16").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("H16").Value = "20000-40000"
# Row 17
$ws.Range("F17").Value = "Systemet (uddyb i bemærkninger)"
$ws.Range("H17").Value = "20000-40000"
# Row 18
$ws.Range("F18").Value = "Strategisk beslutning"
$ws.Range("H18").Value = "20000-40000"
# Row 19
$ws.Range("F19").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("H19").Value = "20000-40000"
# Row 20
$ws.Range("F20").Value = "Ikke oplyst"
$ws.Range("G20").Value = "DataLøn"
$ws.Range("H20").Value = "20000-40000"
# Row 21
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = "DataLøn"
$ws.Range("H21").Value = "20000-40000"
# Row 22
$ws.Range("F22").Value = "Ikke oplyst"
$ws.Range("H22").Value = "20000-40000"
# Row 23
$ws.Range("F23").Value = "Ikke oplyst"
$ws.Range("H23").Value = "20000-40000"
# Row 24
$ws.Range("F24").Value = "Virksomheden lukker"
$ws.Range("H24").Value = "20000-40000"
# Row 25
$ws.Range("F25").Value = "Ikke flere medarbejdere i virksomheden"
$ws.Range("H25").Value = "20000-40000"
# Row 26
$ws.Range("F26").Value = "Ikke oplyst"
$ws.Range("H26").Value = "20000-40000"
# Row 27
$ws.Range("F27").Value = "Utilfredshed (Service - uddyb i bemærkninger)"
$ws.Range("H27").Value = "20000-40000"
# Row 28
$ws.Range("F28").Value = "Virksomheden lukker"
$ws.Range("H28").Value = "20000-40000"
# Row 29
$ws.Range("F29").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("H29").Value = "20000-40000"
# Row 30
$ws.Range("F30").Value = "Systemet (uddyb i bemærkninger)"
$ws.Range("H30").Value = "20000-40000"
# Row 31
$ws.Range("F31").Value = "Strategisk beslutning"
$ws.Range("H31").Value = "20000-40000"
# Row 32
$ws.Range("F32").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("G32").Value = "DataLøn"
$ws.Range("H32").Value = "20000-40000"
# Row 33
$ws.Range("F33").Value = "Ikke oplyst"
$ws.Range("H33").Value = "20000-40000"
# Row 34
$ws.Range("F34").Value = "Ikke oplyst"
$ws.Range("H34").Value = "20000-40000"
# Row 35
$ws.Range("F35").Value = "Ikke oplyst"
$ws.Range("H35").Value = "20000-40000"
# Row 36
$ws.Range("F36").Value = "Bruger ikke produktet"
$ws.Range("H36").Value = "20000-40000"
# Row 37
$ws.Range("F37").Value = "Utilfredshed (Service - uddyb i bemærkninger)"
$ws.Range("H37").Value = "20000-40000"
# Row 38
$ws.Range("F38").Value = "Fusionerer med anden virksomhed"
$ws.Range("H38").Value = "20000-40000"
# Row 39
$ws.Range("F39").Value = "Outsourcing af lønnen (anden leverandør)"
$ws.Range("H39").Value = "20000-40000"
# Row 40
$ws.Range("F40").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("H40").Value = "20000-40000"
